# Contest 28 RCB vs KKR - RCB won.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in the scores for Contest 28 (row 37, "RCB vs KKR")
$ws.Range("E37").Value = 60
$ws.Range("H37").Value = 20
$ws.Range("K37").Value = 80
$ws.Range("N37").Value = 100
$ws.Range("Q37").Value = 40
$ws.Range("T37").Value = 0

# Insert a new row above row 46 (pushes the blank rows down) for Contest 37 (CSK vs RR)
$ws.Rows.Item(46).Insert()

# Re-apply the formatting/formula pattern used by every other contest row,
# copying piecewise so we don't pull in the unused F/I/L/O/R spacer columns.
$ws.Range("A47:E47").Copy()
$ws.Range("A46:E46").PasteSpecial(-4122)

$ws.Range("G47:H47").Copy()
$ws.Range("G46:H46").PasteSpecial(-4122)

$ws.Range("J47:K47").Copy()
$ws.Range("J46:K46").PasteSpecial(-4122)

$ws.Range("M47:N47").Copy()
$ws.Range("M46:N46").PasteSpecial(-4122)

$ws.Range("P47:Q47").Copy()
$ws.Range("P46:Q46").PasteSpecial(-4122)

$ws.Range("S47:T47").Copy()
$ws.Range("S46:T46").PasteSpecial(-4122)

$ws.Range("D46").Formula = "=IF(ISERROR(VLOOKUP(RANK(E46, (`$T46,`$Q46,`$N46,`$K46,`$H46,`$E46), 0),  score, 2, FALSE)),`"`",VLOOKUP(RANK(E46, (`$T46,`$Q46,`$N46,`$K46,`$H46,`$E46), 0),  score, 2, FALSE))"
$ws.Range("G46").Formula = "=IF(ISERROR(VLOOKUP(RANK(H46, (`$T46,`$Q46,`$N46,`$K46,`$H46,`$E46), 0),  score, 2, FALSE)),`"`",VLOOKUP(RANK(H46, (`$T46,`$Q46,`$N46,`$K46,`$H46,`$E46), 0),  score, 2, FALSE))"
$ws.Range("J46").Formula = "=IF(ISERROR(VLOOKUP(RANK(K46, (`$T46,`$Q46,`$N46,`$K46,`$H46,`$E46), 0),  score, 2, FALSE)),`"`",VLOOKUP(RANK(K46, (`$T46,`$Q46,`$N46,`$K46,`$H46,`$E46), 0),  score, 2, FALSE))"
$ws.Range("M46").Formula = "=IF(ISERROR(VLOOKUP(RANK(N46, (`$T46,`$Q46,`$N46,`$K46,`$H46,`$E46), 0),  score, 2, FALSE)),`"`",VLOOKUP(RANK(N46, (`$T46,`$Q46,`$N46,`$K46,`$H46,`$E46), 0),  score, 2, FALSE))"
$ws.Range("P46").Formula = "=IF(ISERROR(VLOOKUP(RANK(Q46, (`$T46,`$Q46,`$N46,`$K46,`$H46,`$E46), 0),  score, 2, FALSE)),`"`",VLOOKUP(RANK(Q46, (`$T46,`$Q46,`$N46,`$K46,`$H46,`$E46), 0),  score, 2, FALSE))"
$ws.Range("S46").Formula = "=IF(ISERROR(VLOOKUP(RANK(T46, (`$T46,`$Q46,`$N46,`$K46,`$H46,`$E46), 0),  score, 2, FALSE)),`"`",VLOOKUP(RANK(T46, (`$T46,`$Q46,`$N46,`$K46,`$H46,`$E46), 0),  score, 2, FALSE))"

$ws.Range("A46").Value = 37
$ws.Range("B46").Value = 1
$ws.Range("C46").Value = "CSK vs RR"

# The "Total" summary row's conditional formatting stayed pinned to row 50;
# re-target it to the row the totals actually live on now (row 51).
$ws.Range("E50").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("E51"))
$ws.Range("H50").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("H51"))
$ws.Range("K50").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("K51"))
$ws.Range("N50").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("N51"))
$ws.Range("Q50").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("Q51"))
$ws.Range("T50").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("T51"))

$ws.Range("U51").Select()
